{"js": "// Rewrite the intranet-PC sentence, then relocate the \"_GoBack\" bookmark\n// from its old spot (an empty trailing paragraph) to the middle of the\n// \"personne de l'ext\u00e9rieur...\" sentence, matching the diff's commit.\n\nconst body = context.document.body;\n\n// 1) \"Pour la partie intranet ... l'intranet.\" sentence gets reworded.\nconst oldSentence =\n  \"Pour la partie intranet il nous est \u00e9galement demand\u00e9 de simuler un PC employ\u00e9 qui doit avoir acc\u00e8s \u00e0 l\\u2019internet et \u00e0 l\\u2019intranet.\";\nconst newSentence =\n  \"Pour la partie intranet il nous est \u00e9galement de\" +\n  \" donner \u00e0 \" +\n  \"un PC employ\u00e9\" +\n  \",\" +\n  \" \" +\n  \"l'\" +\n  \"acc\u00e8s \u00e0 l\\u2019internet et \u00e0 l\\u2019intranet.\";\n\nconst target = body.search(oldSentence, { matchCase: true });\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error(\"Could not locate the intranet sentence to update.\");\n}\ntarget.items[0].insertText(newSentence, \"Replace\");\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark: remove it from its previous (empty\n//    paragraph) location and re-insert it inside the \"ne peut avoir\" run,\n//    right after \"...ne peut av\" (mirrors Word's own last-edit marker).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchorText = body.search(\"personne de l\\u2019ext\u00e9rieur \u00e0 l\\u2019infrastructure ne peut av\", {\n  matchCase: true\n});\nawait context.sync();\n\nif (anchorText.items.length > 0) {\n  const caret = anchorText.items[0].getRange(\"End\");\n  caret.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Reword the \"Pour la partie intranet ... l'intranet.\" sentence and relocate\n# the \"_GoBack\" bookmark from its old (empty trailing paragraph) spot into\n# the middle of the \"personne de l'ext\u00e9rieur...\" sentence - matching the\n# author's commit.\n\n$d = $word.ActiveDocument\n\n# 1) Reword the intranet/PC-employee sentence.\n$oldSentence = \"Pour la partie intranet il nous est \u00e9galement demand\u00e9 de simuler un PC employ\u00e9 qui doit avoir acc\u00e8s \u00e0 l\u2019internet et \u00e0 l\u2019intranet.\"\n$newSentence = \"Pour la partie intranet il nous est \u00e9galement de\" + \" donner \u00e0 \" + \"un PC employ\u00e9\" + \",\" + \" \" + \"l'\" + \"acc\u00e8s \u00e0 l\u2019internet et \u00e0 l\u2019intranet.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($oldSentence)\nif ($found) {\n    # Assigning .Text (rather than Find.Replacement) keeps the literal\n    # straight apostrophe in \"l'\" instead of letting AutoCorrect smarten it.\n    $rng.Text = $newSentence\n}\n\n# 2) Move the \"_GoBack\" bookmark onto the new last-edit location: right\n#    after \"...ne peut av\" in the following sentence.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $d.Content\n$anchorFound = $anchor.Find.Execute(\"personne de l\u2019ext\u00e9rieur \u00e0 l\u2019infrastructure ne peut av\")\nif ($anchorFound) {\n    $anchor.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $anchor)\n}\n"}
